# detection field test data.xlsx - apply field-data update
#  - dog sheet: row 9 "Found" becomes "NA" (was FALSE) + note text tweak
#  - human sheet: wording tweaks on existing notes (rows 2-7) + 4 new rows
#    (Steve x2, Rachael x2) for a 2025-06-06 session

$wb = $excel.ActiveWorkbook

$wsDog = $wb.Worksheets.Item("dog")
$wsHuman = $wb.Worksheets.Item("human")

# ---- dog sheet: row 9 updates ----
$wsDog.Range("H9").Value = "NA"
$wsDog.Range("H9").HorizontalAlignment = -4108   # xlCenter
$wsDog.Range("L9").Value = "Worked downhill. Did not get onto odour. Search got derailed by Koda picking up odour of a frisbee next to the search area. Search abandoned because ignoring a high value item is not part of the controlled evaluation protocol."

# ---- human sheet: reword existing notes to match new "Transects (...)" phrasing ----
$wsHuman.Range("H2").Value = "Transects (long length) with tape marker aids. Became fatigued about half way through and commented the hardest part was staying focused and positive."
$wsHuman.Range("H3").Value = "Transects (short length) with cone marker aids. Fatigue around 3/4 through but better focus. "
$wsHuman.Range("H4").Value = "Transects (short length) no aids. Found rain moth casings. Found it tiring."
$wsHuman.Range("H5").Value = "Transects (short length) no aids. Found GPS! Find was 3/4 of the way through, worked faster than the first attempt."
$wsHuman.Range("H6").Value = "Spiral search from centre out. Very fast find- entirely luck based, GPS was <2m from centre."
$wsHuman.Range("H7").Value = "Spiral search from centre out. More realistic time and search, GPS <10m from centre."

# ---- human sheet: append new rows 8-11 for the 2025-06-06 session ----
# Clone date/time number formats from an existing row so styles are reused.
$wsHuman.Range("B2").Copy()
$wsHuman.Range("B8:B11").PasteSpecial(-4122)   # xlPasteFormats
$wsHuman.Range("C2").Copy()
$wsHuman.Range("C8:C11").PasteSpecial(-4122)   # xlPasteFormats

# Row 8: Steve, spiral search from outside-in then inside-out, did not find
$wsHuman.Range("A8").Value = "Steve"
$wsHuman.Range("B8").Value = 45814
$wsHuman.Range("C8").Value = 0.3888888888888889
$wsHuman.Range("D8").Value = "Overcast, cold"
$wsHuman.Range("E8").Value = $false
$wsHuman.Range("F8").Value = "52 minutes 53 seconds"
$wsHuman.Range("G8").Value = 3173
$wsHuman.Range("H8").Value = "Spiral search from outside in then inside out. "

# Row 9: Steve, transects + perimeter search, did not find (ran out of time)
$wsHuman.Range("A9").Value = "Steve"
$wsHuman.Range("B9").Value = 45814
$wsHuman.Range("C9").Value = 0.44444444444444442
$wsHuman.Range("D9").Value = "Partly sunny, cool."
$wsHuman.Range("E9").Value = $false
$wsHuman.Range("F9").Value = "1 hour"
$wsHuman.Range("G9").Value = 3600
$wsHuman.Range("H9").Value = "Transects (short length) and perimeter search."

# Row 10: Rachael, transects + random walk, found it
$wsHuman.Range("A10").Value = "Rachael"
$wsHuman.Range("B10").Value = 45814
$wsHuman.Range("C10").Value = 0.54513888888888884
$wsHuman.Range("D10").Value = "Overcast, cool."
$wsHuman.Range("E10").Value = $true
$wsHuman.Range("F10").Value = "59 minutes 26 seconds"
$wsHuman.Range("G10").Value = 3566
$wsHuman.Range("H10").Value = "Transects (short length) and random walk."

# Row 11: Rachael, transects + random walk, ran out of time just before finding it
$wsHuman.Range("A11").Value = "Rachael"
$wsHuman.Range("B11").Value = 45814
$wsHuman.Range("C11").Value = 0.60416666666666663
$wsHuman.Range("D11").Value = "Overcast, cool."
$wsHuman.Range("E11").Value = $false
$wsHuman.Range("F11").Value = "1 hour"
$wsHuman.Range("G11").Value = 3600
$wsHuman.Range("H11").Value = "Transects (short length) and random walk. Actually found GPS but after the 1 hour maximum time had elapsed."

# ---- column D (Conditions) is now wider text - refit it ----
$wsHuman.Columns.Item(4).AutoFit()

# ---- sheet view / selection updates ----
# dog sheet scrolls back to A1 and selection moves to L25 (it loses "active" status
# to the human sheet, selected last, below)
$wsDog.Range("A1").Select()
$wsDog.Range("L25").Select()

# human sheet stays the active tab, with the selection parked at B13 (just below
# the newly entered data)
$wsHuman.Range("B13").Select()

Write-Host "Applied detection field test data update."
